# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook

# --- 1. Update time_taken timestamps on the "data" sheet -------------------
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value  = "2021-10-05 14:33:10.610412"
$dataSheet.Range("F3").Value  = "2021-10-05 14:33:10.610419"
$dataSheet.Range("F4").Value  = "2021-10-05 14:33:10.610422"
$dataSheet.Range("F5").Value  = "2021-10-05 14:33:10.610425"
$dataSheet.Range("F6").Value  = "2021-10-05 14:33:10.610428"
$dataSheet.Range("F7").Value  = "2021-10-05 14:33:10.610430"
$dataSheet.Range("F8").Value  = "2021-10-05 14:33:10.610433"
$dataSheet.Range("F9").Value  = "2021-10-05 14:33:10.610435"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:10.610438"
$dataSheet.Range("F11").Value = "2021-10-05 14:33:10.610440"
$dataSheet.Range("F12").Value = "2021-10-05 14:33:10.610443"
$dataSheet.Range("F13").Value = "2021-10-05 14:33:10.610445"
$dataSheet.Range("F14").Value = "2021-10-05 14:33:10.610448"
$dataSheet.Range("F15").Value = "2021-10-05 14:33:10.610450"
$dataSheet.Range("F16").Value = "2021-10-05 14:33:10.610452"

# --- 2. Add a new "metadata" worksheet right after "data" ------------------
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1), bold + centered/top aligned + thin border, like the
# header row on the "data" sheet.
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# A2 uses the same "index" style as column A on the data sheet.
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Font.Bold = $true
$metaSheet.Range("A2").HorizontalAlignment = -4108
$metaSheet.Range("A2").VerticalAlignment = -4160
$metaSheet.Range("A2").Borders.LineStyle = 1

# Data row (row 2).
$metaSheet.Range("B2").Value = "Arrhythmogenic Cardiomyopathy"
$metaSheet.Range("C2").Value = 48

# Keep "0.58" as text, not a number.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.58"

$metaSheet.Range("E2").Value = "2021-05-27T09:08:30.354051Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:10.606524"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/48/?format=json"

# Leave the original "data" sheet active/selected.
$dataSheet.Activate()
